{"js": "// The three \"Poate edita/\u0219terge ...\" list items previously had strikethrough\n// only on the word \"\u0219terge\". The edit extends the strikethrough formatting\n// to the whole line (the complete run of text in each of those paragraphs).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text || \"\";\n  if (text.indexOf(\"Poate edita/\") === 0 && text.indexOf(\"\u0219terge\") !== -1) {\n    paragraph.font.strikeThrough = true;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The three \"Poate edita/\u0219terge ...\" list items previously had strikethrough\n# only on the word \"\u0219terge\". This extends the strikethrough formatting to the\n# entire line (the whole paragraph's run of text).\n\n$doc = $word.ActiveDocument\n\nforeach ($para in $doc.Paragraphs) {\n    $text = $para.Range.Text\n    if ($text -like \"Poate edita/*\" -and $text -like \"*\u0219terge*\") {\n        $para.Range.Font.StrikeThrough = 1\n    }\n}\n"}
